# Updated graphs with drift corrected data
$wb = $excel.ActiveWorkbook

# --- Sheet "Alternative_Input" (sheet1.xml): refresh drift-corrected data table ---
$ws1 = $wb.Worksheets.Item("Alternative_Input")

$ws1.Range("A1").Value = "1631 mNG_FimW FimX_mScI"
$ws1.Range("B1").Value = 20220726
$ws1.Range("C1").Value = "5s interval-2h37"

$ws1.Range("A2").Value = "1631 mNG_FimW FimX_mScI"
$ws1.Range("B2").Value = 20220728
$ws1.Range("C2").Value = "5s interval-2h37"

$ws1.Range("A3").Value = "1631 mNG_FimW FimX_mScI"
$ws1.Range("B3").Value = 20220729
$ws1.Range("C3").Value = "5s interval-2h37"

$ws1.Range("A4").Value = "1631 mNG_FimW FimX_mScI"
$ws1.Range("B4").Value = 20220804
$ws1.Range("C4").Value = "5s interval-2h37"

$ws1.Range("A5").Value = "1632 mNG_FimW FimX_mScI cpdA-"
$ws1.Range("B5").Value = 20220727
$ws1.Range("C5").Value = "5s interval-2h37"

$ws1.Range("A6").Value = "1632 mNG_FimW FimX_mScI cpdA-"
$ws1.Range("B6").Value = 20220728
$ws1.Range("C6").Value = "5s interval-2h37"

$ws1.Range("A7").Value = "1632 mNG_FimW FimX_mScI cpdA-"
$ws1.Range("B7").Value = 20220729
$ws1.Range("C7").Value = "5s interval-2h37"

$ws1.Range("A8").Value = "1633 mNG_FimW FimX_mScI cpdA- pch-"
$ws1.Range("B8").Value = 20220727
$ws1.Range("C8").Value = "5s interval-2h37"

$ws1.Range("A9").Value = "1633 mNG_FimW FimX_mScI cpdA- pch-"
$ws1.Range("B9").Value = 20220728
$ws1.Range("C9").Value = "5s interval-2h37"

$ws1.Range("A10").Value = "1633 mNG_FimW FimX_mScI cpdA- pch-"
$ws1.Range("B10").Value = 20220729
$ws1.Range("C10").Value = "5s interval-2h37"

# --- Sheet "FimW double Jenal Fluo" (sheet2.xml): just refresh the saved selection ---
$ws2 = $wb.Worksheets.Item("FimW double Jenal Fluo")
[void]$ws2.Activate()
[void]$ws2.Range("A46").Select()

# --- Restore the active sheet/selection on "Alternative_Input" ---
[void]$ws1.Activate()
[void]$ws1.Range("A5:XFD6").Select()
